$d = $word.ActiveDocument

# Replace the first paragraph's text.
$d.Content.Find.Execute(
    "-Los requisitos R9 Y R10 son subjetivos.", $false, $false, $false, $false, $false,
    $true, 1, $false, "-  El requisito R9 es subjetivo.", 2)

# Replace the second paragraph's text (keep only one run with the new text).
$d.Content.Find.Execute(
    "-El presupuesto no está justificado. Es una oferta comercial para captar al cliente.",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "- El requisito R10 es subjetivo.", 2)
